$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$NL = [char]10

# ---------------------------------------------------------------------
# Row 9 ("2_4" / Create Buttons on Add word dialog) gets extended:
#  - column C text grows with 3 extra lines
#  - column D flips from "No" to "Yes"
#  - new columns E and F are populated
#  - formatting switches from the plain style (s=1) to the banded
#    fill+border style used by the other data rows (s=2 / s=3)
# ---------------------------------------------------------------------

# Copy the formatting pattern (A=2,B=2,C=3,D=2,E=3,F=3) from row 5, which
# already has exactly that pattern, onto row 9 - this reuses the existing
# style entries instead of creating new ones.
$ws.Range("A5:F5").Copy() | Out-Null
$ws.Range("A9:F9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A9").Value = "2_4"
$ws.Range("B9").Value = "Create Buttons on Add word dialog"
$ws.Range("C9").Value = "Create Add button to save the word into database and cancel button to cancel the work" + $NL + "If use click save the word, database is going to be created if not exists." + $NL + "If database is already exists or created, data will be saved." + $NL + "Database will be ""WORD"" database and ""MEAN"" database to save meaning and word seperatly"
$ws.Range("D9").Value = "Yes"
$ws.Range("E9").Value = "Create SQLite DB Helper class and connect it to app." + $NL + "Get data from EditText and Spinner, and save them into database"
$ws.Range("F9").Value = "MainActivity" + $NL + "DBHelper" + $NL + "AddNewWordDialog" + $NL + "dialog_add_new_word.xml"

$ws.Rows.Item(9).RowHeight = 72

# ---------------------------------------------------------------------
# Brand new row 10 ("2_5" / Add data (word) into ArrayList)
# ---------------------------------------------------------------------

# Same base pattern as row 9 (copied from row 5), then patch E10 back to
# the non-wrapping style (s=2) to match the source formatting exactly.
$ws.Range("A5:F5").Copy() | Out-Null
$ws.Range("A10:F10").PasteSpecial(-4122) | Out-Null
$ws.Range("D5").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A10").Value = "2_5"
$ws.Range("B10").Value = "Add data (word) into ArrayList"
$ws.Range("C10").Value = "The word Database Data needs to be added on ArrayList to be shown" + $NL + "RecyclerView of MainActivity."
$ws.Range("D10").Value = "Yes"
$ws.Range("E10").Value = "ArrayList that contains WordsList connected to Adapter to show on Recyclerview."
$ws.Range("F10").Value = "MainActivity" + $NL + "MainActivityRecyclerViewAdapter" + $NL + "DBHelper" + $NL + "AddNewWordDialog" + $NL + "WordsList" + $NL + "dialog_add_new_word.xml" + $NL + "rvmainwords_items.xml"

$ws.Rows.Item(10).RowHeight = 115.2

# ---------------------------------------------------------------------
# Column E gets wider to fit the new, longer text
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 67.83333333333333

# ---------------------------------------------------------------------
# Scroll the view down a couple of rows (topLeftCell A4 -> A6) while
# keeping the existing C7 selection
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
